# Add season-record columns (Wins / Losses / Ties) to the SFG 2010 sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): copy the formatting of an existing header cell
# (bold, centered, bordered) onto the three new header cells, then set
# their text.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# --- Data rows (2-44): every player on the roster shares the same team
# season record: 92 wins, 70 losses, 0 ties.
$lastRow = 44
for ($row = 2; $row -le $lastRow; $row++) {
    $ws.Cells.Item($row, 30).Value = 92
    $ws.Cells.Item($row, 31).Value = 70
    $ws.Cells.Item($row, 32).Value = 0
}
